# SCD0013-004 - Admin SLN melakukan Modul Mapping
# Rename sheet from SCD0219 -> SCD0013 and update the TC_ID value from
# "DGS-234" to "SCD0013-004" (column B, row 2). Also reflect the resulting
# column B autosize and the new active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "SCD0013"

# Update the TC_ID cell with the new test-case id
$ws.Range("B2").Value = "SCD0013-004"

# Widen column B so the longer TC_ID value fits (closest attainable width)
$ws.Columns("B").ColumnWidth = 11.666667

# Move the active selection to B3, matching the post-edit cursor position
$ws.Range("B3").Select()
